$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(44503, 13146, 2539, 15685),
    @(44504, 13166, 2545, 15711),
    @(44505, 13205, 2551, 15756),
    @(44506, 13229, 2551, 15780),
    @(44507, 13269, 2551, 15820),
    @(44508, 13314, 2556, 15870),
    @(44509, 13367, 2571, 15938),
    @(44510, 13405, 2580, 15985),
    @(44511, 13446, 2584, 16030),
    @(44512, 13485, 2589, 16074),
    @(44513, 13537, 2594, 16131),
    @(44514, 13598, 2596, 16194),
    @(44515, 13644, 2602, 16246),
    @(44516, 13687, 2615, 16302),
    @(44517, 13725, 2616, 16341)
)

$startRow = 384
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
